# Applies numeric value updates to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the Hyperion_Profits scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 3815.5
$ws.Range("I11").Value = 3815.5
$ws.Range("K11").Value = 3815.5
$ws.Range("M11").Value = -3675.5

# Row 33
$ws.Range("H33").Value = 496.75
$ws.Range("I33").Value = 608
$ws.Range("J33").Value = 290.14285
$ws.Range("K33").Value = 608
$ws.Range("L33").Value = 290.14285
$ws.Range("M33").Value = -379
$ws.Range("N33").Value = -748.14285

# Row 40
$ws.Range("H40").Value = 4409.524
$ws.Range("I40").Value = 4131.4375
$ws.Range("J40").Value = 5299.4
$ws.Range("K40").Value = 4131.4375
$ws.Range("L40").Value = 5299.4
$ws.Range("M40").Value = -3956.4375
$ws.Range("N40").Value = -5649.4

# Row 107
$ws.Range("H107").Value = 151.5
$ws.Range("J107").Value = 243.57143
$ws.Range("L107").Value = 243.57143
$ws.Range("N107").Value = -4083.57143

# Row 111
$ws.Range("H111").Value = 2074.5715
$ws.Range("I111").Value = 2390.3076
$ws.Range("J111").Value = 1561.5
$ws.Range("K111").Value = 7170.9228
$ws.Range("L111").Value = 4684.5
$ws.Range("M111").Value = -4103.9228
$ws.Range("N111").Value = -10818.5

# Row 112
$ws.Range("H112").Value = 4478.9443
$ws.Range("J112").Value = 5088.968
$ws.Range("L112").Value = 15266.904
$ws.Range("N112").Value = -17482.904

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 47150
$ws.Range("I2").Value = 57091.668
$ws.Range("J2").Value = 2412.5
$ws.Range("K2").Value = 57091.668
$ws.Range("L2").Value = 2412.5
$ws.Range("M2").Value = -56978.668
$ws.Range("N2").Value = -2638.5

# Row 45
$ws.Range("H45").Value = 76738.92999999999
$ws.Range("I45").Value = 127045.375
$ws.Range("J45").Value = 9663.666999999999
$ws.Range("K45").Value = 127045.375
$ws.Range("L45").Value = 9663.666999999999
$ws.Range("M45").Value = -126668.375
$ws.Range("N45").Value = -10417.667

# Row 97
$ws.Range("H97").Value = 9130.272000000001
$ws.Range("J97").Value = 2439.8333
$ws.Range("L97").Value = 2439.8333
$ws.Range("N97").Value = -3431.8333

# Row 102
$ws.Range("H102").Value = 5138.375
$ws.Range("I102").Value = 5138.375
$ws.Range("K102").Value = 5138.375
$ws.Range("M102").Value = -3516.375

# Row 110
$ws.Range("H110").Value = 5286.885
$ws.Range("I110").Value = 1053.84
$ws.Range("K110").Value = 1053.84
$ws.Range("M110").Value = 991.1600000000001

# Row 116
$ws.Range("H116").Value = 47150
$ws.Range("I116").Value = 57091.668
$ws.Range("J116").Value = 2412.5
$ws.Range("K116").Value = 57091.668
$ws.Range("L116").Value = 2412.5
$ws.Range("M116").Value = -54797.668
$ws.Range("N116").Value = -7000.5

# Row 119
$ws.Range("H119").Value = 53492.645
$ws.Range("J119").Value = 53492.645
$ws.Range("L119").Value = 53492.645
$ws.Range("N119").Value = -63168.645

# Row 132
$ws.Range("H132").Value = 3224.08
$ws.Range("I132").Value = 2305.158
$ws.Range("J132").Value = 6134
$ws.Range("K132").Value = 6915.474
$ws.Range("L132").Value = 18402
$ws.Range("M132").Value = -4385.474
$ws.Range("N132").Value = -23462

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 47150
$ws.Range("I3").Value = 57091.668
$ws.Range("J3").Value = 2412.5
$ws.Range("K3").Value = 57091.668
$ws.Range("L3").Value = 2412.5
$ws.Range("M3").Value = -56977.668
$ws.Range("N3").Value = -2640.5

# Row 86
$ws.Range("H86").Value = 5398.4644
$ws.Range("I86").Value = 7035.1577
$ws.Range("J86").Value = 1943.2222
$ws.Range("K86").Value = 7035.1577
$ws.Range("L86").Value = 1943.2222
$ws.Range("M86").Value = -5912.1577
$ws.Range("N86").Value = -4189.2222

# Row 89
$ws.Range("H89").Value = 5398.4644
$ws.Range("I89").Value = 7035.1577
$ws.Range("J89").Value = 1943.2222
$ws.Range("K89").Value = 35175.7885
$ws.Range("L89").Value = 9716.110999999999
$ws.Range("M89").Value = -29559.7885
$ws.Range("N89").Value = -20948.111

# Row 95
$ws.Range("H95").Value = 9809.5
$ws.Range("J95").Value = 9809.5
$ws.Range("L95").Value = 9809.5
$ws.Range("N95").Value = -15301.5

# Row 99
$ws.Range("H99").Value = 3029.6316
$ws.Range("I99").Value = 2429.923
$ws.Range("J99").Value = 4329
$ws.Range("K99").Value = 2429.923
$ws.Range("L99").Value = 4329
$ws.Range("M99").Value = -931.9229999999998
$ws.Range("N99").Value = -7325

# Row 107
$ws.Range("H107").Value = 2332.36
$ws.Range("I107").Value = 2149.1052
$ws.Range("K107").Value = 2149.1052
$ws.Range("M107").Value = -229.1052

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1498.8235
$ws.Range("I16").Value = 1387.7333
$ws.Range("J16").Value = 2332
$ws.Range("K16").Value = 1387.7333
$ws.Range("L16").Value = 2332
$ws.Range("M16").Value = -1100.7333
$ws.Range("N16").Value = -2906

# Row 26
$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -713

# Row 76
$ws.Range("H76").Value = 5166.6665
$ws.Range("I76").Value = 5166.6665
$ws.Range("K76").Value = 5166.6665
$ws.Range("M76").Value = -4851.6665

# Row 79
$ws.Range("H79").Value = 5166.6665
$ws.Range("I79").Value = 5166.6665
$ws.Range("K79").Value = 5166.6665
$ws.Range("M79").Value = -4074.6665

# Row 86
$ws.Range("H86").Value = 11483.9
$ws.Range("I86").Value = 6815.3335
$ws.Range("K86").Value = 6815.3335
$ws.Range("M86").Value = -5692.3335

# Row 89
$ws.Range("H89").Value = 11483.9
$ws.Range("I89").Value = 6815.3335
$ws.Range("K89").Value = 34076.6675
$ws.Range("M89").Value = -28460.6675

# Row 99
$ws.Range("H99").Value = 4061.1765
$ws.Range("I99").Value = 3528.8333
$ws.Range("K99").Value = 3528.8333
$ws.Range("M99").Value = -2030.8333

# Row 108
$ws.Range("H108").Value = 43991
$ws.Range("J108").Value = 52488.75
$ws.Range("L108").Value = 52488.75
$ws.Range("N108").Value = -60168.75

# Row 113
$ws.Range("H113").Value = 1498.8235
$ws.Range("I113").Value = 1387.7333
$ws.Range("J113").Value = 2332
$ws.Range("K113").Value = 1387.7333
$ws.Range("L113").Value = 2332
$ws.Range("M113").Value = 782.2666999999999
$ws.Range("N113").Value = -6672

# Row 126
$ws.Range("H126").Value = 4061.1765
$ws.Range("I126").Value = 3528.8333
$ws.Range("K126").Value = 10586.4999
$ws.Range("M126").Value = -8116.499899999999

# Row 134
$ws.Range("H134").Value = 34415.965
$ws.Range("I134").Value = 47173.55
$ws.Range("K134").Value = 141520.65
$ws.Range("M134").Value = -138985.65

# Row 141
$ws.Range("H141").Value = 132617.56
$ws.Range("J141").Value = 132617.56
$ws.Range("L141").Value = 132617.56
$ws.Range("N141").Value = -142977.56

$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 1140.3636
$ws.Range("I50").Value = 864
$ws.Range("K50").Value = 2592
$ws.Range("M50").Value = -2111

# Row 53
$ws.Range("H53").Value = 1140.3636
$ws.Range("I53").Value = 864
$ws.Range("K53").Value = 2592
$ws.Range("M53").Value = -2111

# Row 63
$ws.Range("H63").Value = 10092.5
$ws.Range("I63").Value = 2987.5
$ws.Range("J63").Value = 11868.75
$ws.Range("K63").Value = 8962.5
$ws.Range("L63").Value = 35606.25
$ws.Range("M63").Value = -8213.5
$ws.Range("N63").Value = -37104.25

# Row 66
$ws.Range("H66").Value = 10092.5
$ws.Range("I66").Value = 2987.5
$ws.Range("J66").Value = 11868.75
$ws.Range("K66").Value = 26887.5
$ws.Range("L66").Value = 106818.75
$ws.Range("M66").Value = -23143.5
$ws.Range("N66").Value = -114306.75

# Row 92
$ws.Range("H92").Value = 1218.4
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -17496

# Row 104
$ws.Range("H104").Value = 2500
$ws.Range("J104").Value = 2500
$ws.Range("L104").Value = 7500
$ws.Range("N104").Value = -12742

# Row 132
$ws.Range("H132").Value = 1857.875
$ws.Range("J132").Value = 2177.7778
$ws.Range("L132").Value = 19600.0002
$ws.Range("N132").Value = -24660.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 3412.6365
$ws.Range("I113").Value = 1931
$ws.Range("K113").Value = 1931
$ws.Range("M113").Value = 239

# Row 117
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884

# Row 122
$ws.Range("H122").Value = 426750.1
$ws.Range("I122").Value = 810440.75
$ws.Range("J122").Value = 4690.4
$ws.Range("K122").Value = 2431322.25
$ws.Range("L122").Value = 14071.2
$ws.Range("M122").Value = -2428872.25
$ws.Range("N122").Value = -18971.2

$ws = $wb.Worksheets.Item("LTW")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# Row 82
$ws.Range("H82").Value = 125002560
$ws.Range("I82").Value = 200002100
$ws.Range("J82").Value = 3333.3333
$ws.Range("K82").Value = 200002100
$ws.Range("L82").Value = 3333.3333
$ws.Range("M82").Value = -200001739
$ws.Range("N82").Value = -4055.3333

# Row 85
$ws.Range("H85").Value = 125002560
$ws.Range("I85").Value = 200002100
$ws.Range("J85").Value = 3333.3333
$ws.Range("K85").Value = 200002100
$ws.Range("L85").Value = 3333.3333
$ws.Range("M85").Value = -200000852
$ws.Range("N85").Value = -5829.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 853.63336
$ws.Range("I113").Value = 440.125
$ws.Range("J113").Value = 1004
$ws.Range("K113").Value = 1320.375
$ws.Range("L113").Value = 3012
$ws.Range("M113").Value = 849.625
$ws.Range("N113").Value = -7352

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# Row 126
$ws.Range("H126").Value = 3062.875
$ws.Range("J126").Value = 2299.6667
$ws.Range("L126").Value = 6899.000100000001
$ws.Range("N126").Value = -11839.0001
